# Update the LR-pairs data (Efnb3-Ephb1) with the newly computed TPM-based
# NATMI statistics. The sending/target cluster set grows from
# {Inflammatory-Mac, MuSCs, Neutrophils, Resolving-Mac} x {ECs, MuSCs} (8 rows)
# to {ECs, Inflammatory-Mac, MuSCs, Resolving-Mac} x {ECs, MuSCs, Neutrophils}
# (12 rows), so rows 2-13 are (re)written in full and the sheet grows from
# A1:T9 to A1:T13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Efnb3", "Ephb1", "ECs", 1, 0.3333333333333333, 0.2163973333333333, 0.649192, 0.2690509661419307, 0.2690509661419307, 3, 1, 1.397474666666667, 4.192424, 0.9253320219903044, 0.9253320219903044, 0.3024097912675556, 2.721688121408, 0.2489614745185577, 0.2489614745185577),
    @("ECs", "Efnb3", "Ephb1", "MuSCs", 1, 0.3333333333333333, 0.2163973333333333, 0.649192, 0.2690509661419307, 0.2690509661419307, 3, 1, 0.07803566666666666, 0.234107, 0.05167099121464913, 0.05167099121464913, 0.01688671017155555, 0.151980391544, 0.01390213010781256, 0.01390213010781256),
    @("ECs", "Efnb3", "Ephb1", "Neutrophils", 1, 0.3333333333333333, 0.2163973333333333, 0.649192, 0.2690509661419307, 0.2690509661419307, 1, 0.3333333333333333, 0.034731, 0.104193, 0.02299698679504644, 0.02299698679504644, 0.007515695784, 0.06764126205599999, 0.006187361515560467, 0.006187361515560466),
    @("Inflammatory-Mac", "Efnb3", "Ephb1", "ECs", 1, 0.3333333333333333, 0.06315533333333333, 0.189466, 0.07852224049441002, 0.07852224049441003, 3, 1, 1.397474666666667, 4.192424, 0.9253320219903044, 0.9253320219903044, 0.08825797839822221, 0.794321805584, 0.07265914356790139, 0.0726591435679014),
    @("Inflammatory-Mac", "Efnb3", "Ephb1", "MuSCs", 1, 0.3333333333333333, 0.06315533333333333, 0.189466, 0.07852224049441002, 0.07852224049441003, 3, 1, 0.07803566666666666, 0.234107, 0.05167099121464913, 0.05167099121464913, 0.004928368540222221, 0.044355316862, 0.004057321998741226, 0.004057321998741227),
    @("Inflammatory-Mac", "Efnb3", "Ephb1", "Neutrophils", 1, 0.3333333333333333, 0.06315533333333333, 0.189466, 0.07852224049441002, 0.07852224049441003, 1, 0.3333333333333333, 0.034731, 0.104193, 0.02299698679504644, 0.02299698679504644, 0.002193447882, 0.019741030938, 0.001805774927767408, 0.001805774927767408),
    @("MuSCs", "Efnb3", "Ephb1", "ECs", 2, 0.6666666666666666, 0.4426526666666666, 1.327958, 0.5503585732663155, 0.5503585732663157, 3, 1, 1.397474666666667, 4.192424, 0.9253320219903044, 0.9253320219903044, 0.618595887799111, 5.567362990192, 0.5092644114202188, 0.5092644114202189),
    @("MuSCs", "Efnb3", "Ephb1", "MuSCs", 2, 0.6666666666666666, 0.4426526666666666, 1.327958, 0.5503585732663155, 0.5503585732663157, 3, 1, 0.07803566666666666, 0.234107, 0.05167099121464913, 0.05167099121464913, 0.03454269594511111, 0.310884263506, 0.02843757300415062, 0.02843757300415063),
    @("MuSCs", "Efnb3", "Ephb1", "Neutrophils", 2, 0.6666666666666666, 0.4426526666666666, 1.327958, 0.5503585732663155, 0.5503585732663157, 1, 0.3333333333333333, 0.034731, 0.104193, 0.02299698679504644, 0.02299698679504644, 0.015373769766, 0.138363927894, 0.01265658884194606, 0.01265658884194606),
    @("Resolving-Mac", "Efnb3", "Ephb1", "ECs", 1, 0.3333333333333333, 0.08209333333333334, 0.24628, 0.1020682200973436, 0.1020682200973436, 3, 1, 1.397474666666667, 4.192424, 0.9253320219903044, 0.9253320219903044, 0.1147233536355556, 1.03251018272, 0.09444699248362638, 0.09444699248362638),
    @("Resolving-Mac", "Efnb3", "Ephb1", "MuSCs", 1, 0.3333333333333333, 0.08209333333333334, 0.24628, 0.1020682200973436, 0.1020682200973436, 3, 1, 0.07803566666666666, 0.234107, 0.05167099121464913, 0.05167099121464913, 0.006406207995555555, 0.05765587195999999, 0.005273966103944715, 0.005273966103944715),
    @("Resolving-Mac", "Efnb3", "Ephb1", "Neutrophils", 1, 0.3333333333333333, 0.08209333333333334, 0.24628, 0.1020682200973436, 0.1020682200973436, 1, 0.3333333333333333, 0.034731, 0.104193, 0.02299698679504644, 0.02299698679504644, 0.00285118356, 0.02566065204, 0.002347261509772505, 0.002347261509772504)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $rowVals = $data[$i]
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}
